$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("G2").Value = 49
$ws.Range("I2").Value = 43
$ws.Range("C3").Value = 46
$ws.Range("E3").Value = 62
$ws.Range("E6").Value = 4
$ws.Range("C9").Value = 217
$ws.Range("D9").Value = 192
$ws.Range("E9").Value = 203
$ws.Range("F9").Value = 237
$ws.Range("G9").Value = 220
$ws.Range("H9").Value = 194
$ws.Range("I9").Value = 259
$ws.Range("J9").Value = 184
$ws.Range("B10").Value = 502
$ws.Range("C10").Value = 614
$ws.Range("D10").Value = 808
$ws.Range("E10").Value = 913
$ws.Range("F10").Value = 1026
$ws.Range("G10").Value = 545
$ws.Range("H10").Value = 212
$ws.Range("I10").Value = 383
$ws.Range("J10").Value = 335
$ws.Range("B11").Value = 736
$ws.Range("C11").Value = 910
$ws.Range("D11").Value = 1110
$ws.Range("E11").Value = 1217
$ws.Range("F11").Value = 1364
$ws.Range("G11").Value = 873
$ws.Range("H11").Value = 513
$ws.Range("I11").Value = 788
$ws.Range("J11").Value = 677

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("G2").Value = 3
$ws.Range("I2").Value = 3
$ws.Range("E5").Value = 1
$ws.Range("F6").Value = 25
$ws.Range("G6").Value = 28
$ws.Range("E7").Value = 36
$ws.Range("J7").Value = 18
$ws.Range("E8").Value = 56
$ws.Range("F8").Value = 83
$ws.Range("G8").Value = 62
$ws.Range("I8").Value = 43
$ws.Range("J8").Value = 32

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("G7").Value = 9
$ws.Range("G8").Value = 18

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("G6").Value = 4
$ws.Range("I6").Value = 11
$ws.Range("J6").Value = 9
$ws.Range("G8").Value = 19
$ws.Range("I8").Value = 39
$ws.Range("J8").Value = 25

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("B8").Value = 64
$ws.Range("C8").Value = 114
$ws.Range("E8").Value = 265
$ws.Range("F8").Value = 303
$ws.Range("I8").Value = 96
$ws.Range("B9").Value = 85
$ws.Range("C9").Value = 139
$ws.Range("E9").Value = 310
$ws.Range("F9").Value = 338
$ws.Range("I9").Value = 158

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 19
$ws.Range("D7").Value = 24

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 16

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("F6").Value = 14
$ws.Range("F8").Value = 32

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("C8").Value = 44
$ws.Range("G19").Value = 18
$ws.Range("G23").Value = 3
$ws.Range("C28").Value = 66
$ws.Range("E28").Value = 47
$ws.Range("J28").Value = 27
$ws.Range("E29").Value = 12
$ws.Range("E32").Value = 56
$ws.Range("F32").Value = 83
$ws.Range("G32").Value = 62
$ws.Range("I32").Value = 43
$ws.Range("J32").Value = 32
$ws.Range("G36").Value = 19
$ws.Range("I36").Value = 39
$ws.Range("J36").Value = 25
$ws.Range("J38").Value = 2
$ws.Range("B47").Value = 15
$ws.Range("E47").Value = 26
$ws.Range("D48").Value = 11
$ws.Range("D50").Value = 16
$ws.Range("D52").Value = 16
$ws.Range("B53").Value = 85
$ws.Range("C53").Value = 139
$ws.Range("E53").Value = 310
$ws.Range("F53").Value = 338
$ws.Range("I53").Value = 158
$ws.Range("J61").Value = 11
$ws.Range("F65").Value = 32
$ws.Range("D70").Value = 24
$ws.Range("D74").Value = 44
$ws.Range("F76").Value = 23
$ws.Range("H76").Value = 27
$ws.Range("J76").Value = 25
$ws.Range("B77").Value = 30
$ws.Range("D78").Value = 28
$ws.Range("H78").Value = 8
$ws.Range("J79").Value = 5
$ws.Range("F82").Value = 9
$ws.Range("D87").Value = 6
$ws.Range("E87").Value = 9
$ws.Range("H94").Value = 8
$ws.Range("I94").Value = 9
$ws.Range("B98").Value = 736
$ws.Range("C98").Value = 910
$ws.Range("D98").Value = 1110
$ws.Range("E98").Value = 1217
$ws.Range("F98").Value = 1364
$ws.Range("G98").Value = 873
$ws.Range("H98").Value = 513
$ws.Range("I98").Value = 788
$ws.Range("J98").Value = 677

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("J6").Value = 2
$ws.Range("J7").Value = 5

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("H4").Value = 4
$ws.Range("D5").Value = 26
$ws.Range("D6").Value = 28
$ws.Range("H6").Value = 8

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("C3").Value = 2
$ws.Range("C6").Value = 20
$ws.Range("E6").Value = 11
$ws.Range("C7").Value = 44
$ws.Range("J7").Value = 14
$ws.Range("C8").Value = 66
$ws.Range("E8").Value = 47
$ws.Range("J8").Value = 27

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("E3").Value = 1
$ws.Range("B7").Value = 15
$ws.Range("B8").Value = 15
$ws.Range("E8").Value = 26

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 9
$ws.Range("E7").Value = 12

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("H6").Value = 9
$ws.Range("F7").Value = 14
$ws.Range("H7").Value = 16
$ws.Range("J7").Value = 19
$ws.Range("F8").Value = 23
$ws.Range("H8").Value = 27
$ws.Range("J8").Value = 25

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = 9

$ws = $wb.Worksheets.Item("River North")
$ws.Range("D6").Value = 37
$ws.Range("D7").Value = 44

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("D6").Value = 12
$ws.Range("D7").Value = 16

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 3

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("B8").Value = 17
$ws.Range("B9").Value = 30

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 3
$ws.Range("H7").Value = 8
$ws.Range("I7").Value = 9

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 11

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("J5").Value = 2
$ws.Range("J6").Value = 2

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("C7").Value = 21
$ws.Range("C8").Value = 44

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 4
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 9
